$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The commit removes worker JONATHAN ALFREDO ORTEGA TAPIA (row 17) from the
# "Estado de Cuenta" worker-mora table, leaving only LUCAS JOSE PATERNINA TAPIA
# (row 16). Deleting the whole row shifts the subsequent signature rows
# (22/23 -> 21/22) and shrinks the merged cells / used range automatically.
$ws.Rows("17:17").Delete()

# Update the aggregate "VALOR MORA" total (E11) to reflect only the
# remaining worker's mora value (previously the sum of both workers).
$ws.Range("E11").Value = 58667

# Update worker / period counts (previously 2 workers / 2 periods, now 1/1).
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Column D ("Nombre Trabajador") auto-fit width shrinks now that the longest
# name in the table ("JONATHAN ALFREDO ORTEGA TAPIA") was removed.
$ws.Columns.Item(4).ColumnWidth = 27.65
